$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rushing": two new players (T.Homer, D.Metcalf) were logged for
# Week 17, inserted into the existing table (not appended), so every row
# below each insertion point shifts down by one.
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Insert a row for T.Homer right after A.Collins (row 4), before D.Dallas.
$rushing.Rows.Item(5).Insert()
$rushing.Range("A4").Copy()
$rushing.Range("A5").PasteSpecial(-4122)

# Insert a row for D.Metcalf right after A.Peterson (now row 7), before T.Lockett.
$rushing.Rows.Item(8).Insert()
$rushing.Range("A7").Copy()
$rushing.Range("A8").PasteSpecial(-4122)

# Now write every data row (1DATT, 2DATT, 3DATT, RZATT) for the final, post-insert
# layout of the Rushing sheet.
$rushingRows = @(
    @{ Row = 2;  Num = 0;  Name = "R.Wilson";   C = 10; D = 11; E = 15; F = 4 }
    @{ Row = 3;  Num = 1;  Name = "R.Penny";    C = 55; D = 30; E = 5;  F = 14 }
    @{ Row = 4;  Num = 2;  Name = "A.Collins";  C = 55; D = 36; E = 9;  F = 12 }
    @{ Row = 5;  Num = 3;  Name = "T.Homer";    C = 3;  D = 3;  E = 0;  F = 1 }
    @{ Row = 6;  Num = 4;  Name = "D.Dallas";   C = 4;  D = 3;  E = 4;  F = 2 }
    @{ Row = 7;  Num = 5;  Name = "A.Peterson"; C = 7;  D = 2;  E = 1;  F = 5 }
    @{ Row = 8;  Num = 6;  Name = "D.Metcalf";  C = 1;  D = 0;  E = 0;  F = 0 }
    @{ Row = 9;  Num = 7;  Name = "T.Lockett";  C = 2;  D = 0;  E = 0;  F = 0 }
    @{ Row = 10; Num = 8;  Name = "F.Swain";    C = 3;  D = 2;  E = 0;  F = 0 }
    @{ Row = 11; Num = 9;  Name = "D.Eskridge"; C = 4;  D = 1;  E = 0;  F = 0 }
    @{ Row = 12; Num = 10; Name = "G.Everett";  C = 2;  D = 1;  E = 0;  F = 0 }
    @{ Row = 13; Num = 11; Name = "W.Dissly";   C = 0;  D = 1;  E = 0;  F = 0 }
)

foreach ($r in $rushingRows) {
    $row = $r.Row
    $rushing.Range("A$row").Value = $r.Num
    $rushing.Range("B$row").Value = $r.Name
    $rushing.Range("C$row").Value = $r.C
    $rushing.Range("D$row").Value = $r.D
    $rushing.Range("E$row").Value = $r.E
    $rushing.Range("F$row").Value = $r.F
}

# ---------------------------------------------------------------------------
# Sheet "Receiving": same 12 players as before (no rows added/removed),
# just Week 17 stats logged on top of the running totals.
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

$receivingRows = @(
    @{ Row = 2;  Num = 0;  Name = "R.Penny";     C = 14; D = 11; E = 1;  F = 0;  G = 2;  H = 1 }
    @{ Row = 3;  Num = 1;  Name = "A.Collins";   C = 16; D = 14; E = 1;  F = 1;  G = 0;  H = 0 }
    @{ Row = 4;  Num = 2;  Name = "D.Dallas";    C = 13; D = 11; E = 0;  F = 0;  G = 4;  H = 3 }
    @{ Row = 5;  Num = 3;  Name = "N.Bellore";   C = 1;  D = 1;  E = 0;  F = 0;  G = 0;  H = 0 }
    @{ Row = 6;  Num = 4;  Name = "D.Metcalf";   C = 84; D = 60; E = 34; F = 10; G = 19; H = 13 }
    @{ Row = 7;  Num = 5;  Name = "T.Lockett";   C = 67; D = 51; E = 44; F = 25; G = 8;  H = 4 }
    @{ Row = 8;  Num = 6;  Name = "F.Swain";     C = 31; D = 19; E = 8;  F = 5;  G = 4;  H = 1 }
    @{ Row = 9;  Num = 7;  Name = "D.Eskridge";  C = 13; D = 10; E = 7;  F = 0;  G = 4;  H = 2 }
    @{ Row = 10; Num = 8;  Name = "P.Hart";      C = 10; D = 7;  E = 2;  F = 0;  G = 1;  H = 0 }
    @{ Row = 11; Num = 9;  Name = "G.Everett";   C = 55; D = 43; E = 5;  F = 4;  G = 7;  H = 4 }
    @{ Row = 12; Num = 10; Name = "W.Dissly";    C = 21; D = 16; E = 5;  F = 5;  G = 3;  H = 2 }
    @{ Row = 13; Num = 11; Name = "C.Parkinson"; C = 4;  D = 3;  E = 1;  F = 0;  G = 1;  H = 0 }
)

foreach ($r in $receivingRows) {
    $row = $r.Row
    $receiving.Range("A$row").Value = $r.Num
    $receiving.Range("B$row").Value = $r.Name
    $receiving.Range("C$row").Value = $r.C
    $receiving.Range("D$row").Value = $r.D
    $receiving.Range("E$row").Value = $r.E
    $receiving.Range("F$row").Value = $r.F
    $receiving.Range("G$row").Value = $r.G
    $receiving.Range("H$row").Value = $r.H
}

Write-Host "Applied Week 17 data to Rushing and Receiving sheets."
